$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data as plain text values,
# using a leading apostrophe to force text interpretation (avoids Excel
# auto-converting numeric-looking strings to numbers), then resetting the
# cell style back to Normal so no stray number-format/style is introduced.
$ws.Range("D2").Value = "'59.037.16"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.58%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.498.60"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.54%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.02%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'537.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.38%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'137.18"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -1.68%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.11%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +1.58%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.520.54"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.26%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -0.61%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -2.08%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'5.33"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -1.00%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  -3.11%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'2.942.26"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.74%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'23.03"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -1.61%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'58.887.00"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.66%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  -1.43%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.511.62"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.12%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'11.13"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.26%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -0.61%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'323.76"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.31%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -0.02%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +2.06%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'65.90"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +4.23%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.422"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.31%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -1.68%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -0.44%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'7.55"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -3.51%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'6.70"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -2.38%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.0₃0772"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.81%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -1.37%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'167.76"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +2.07%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +4.89%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -0.16%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  +0.68%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -0.31%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -3.57%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -3.40%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'36.69"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.70%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.812"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +0.26%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'3.62"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -1.69%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'284.44"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +1.07%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'5.14"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -1.78%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'132.80"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +6.87%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -0.16%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  +1.65%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'10.90"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.49%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.0926"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -1.04%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -1.28%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -2.13%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'17.35"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -2.90%  "
$ws.Range("E51").Style = "Normal"
